$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. First paragraph: "This is a Microsoft word document." gets two trailing
#    spaces appended (same run), followed by three additional runs colored
#    C00000 (dark red) that spell out the parenthetical remark.
# ---------------------------------------------------------------------------

function Append-ColoredText($ParagraphIndex, $Text, $Color) {
    $p = $d.Paragraphs.Item($ParagraphIndex)
    $startPos = $p.Range.End - 1
    $ins = $d.Range($startPos, $startPos)
    $ins.InsertAfter($Text)
    $endPos = $startPos + $Text.Length
    $applied = $d.Range($startPos, $endPos)
    $applied.Font.Color = $Color
}

$p1 = $d.Paragraphs.Item(1)
$endOfP1 = $p1.Range.End - 1
$plain = $d.Range($endOfP1, $endOfP1)
$plain.InsertAfter("  ")

$darkRed = 192  # 0x0000C0 (BGR) == w:val="C00000"

Append-ColoredText 1 "(This is a change – Ve" $darkRed
Append-ColoredText 1 "rsion for branch alternate" $darkRed
Append-ColoredText 1 ")" $darkRed

# ---------------------------------------------------------------------------
# 2. Insert a new, empty paragraph right after the "It will be treated..."
#    paragraph. It carries paragraph-level shading and paragraph-mark
#    run-properties only (no visible run), matching a "press Enter then set
#    formatting without typing" edit.
# ---------------------------------------------------------------------------

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(3)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="202122"/></w:rPr></w:pPr></w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($xml)
